$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '28.296.69'
Set-TextValue $ws.Range("E2") '  -0.73%  '
Set-TextValue $ws.Range("D3") '1.810.50'
Set-TextValue $ws.Range("E3") '  -0.87%  '
Set-TextValue $ws.Range("D4") '0.9992'
Set-TextValue $ws.Range("E4") '  -0.54%  '
Set-TextValue $ws.Range("D5") '312.83'
Set-TextValue $ws.Range("E5") '  -1.56%  '
Set-TextValue $ws.Range("D6") '0.9990'
Set-TextValue $ws.Range("E6") '  -0.55%  '
Set-TextValue $ws.Range("D7") '0.5163'
Set-TextValue $ws.Range("E7") '  +0.01%  '
Set-TextValue $ws.Range("D8") '0.4011'
Set-TextValue $ws.Range("E8") '  +3.86%  '
Set-TextValue $ws.Range("D9") '0.07866'
Set-TextValue $ws.Range("E9") '  -5.54%  '
Set-TextValue $ws.Range("D10") '1.112'
Set-TextValue $ws.Range("E10") '  -0.54%  '
Set-TextValue $ws.Range("D11") '41.13'
Set-TextValue $ws.Range("E11") '  -2.12%  '
Set-TextValue $ws.Range("D12") '6.362'
Set-TextValue $ws.Range("E12") '  -0.77%  '
Set-TextValue $ws.Range("D13") '0.9988'
Set-TextValue $ws.Range("E13") '  -0.54%  '
Set-TextValue $ws.Range("D14") '20.47'
Set-TextValue $ws.Range("E14") '  -3.19%  '
Set-TextValue $ws.Range("D15") '7.346'
Set-TextValue $ws.Range("E15") '  -1.93%  '
Set-TextValue $ws.Range("D16") '1.799.32'
Set-TextValue $ws.Range("D17") '93.00'
Set-TextValue $ws.Range("E17") '  -1.11%  '
Set-TextValue $ws.Range("D18") '0.00001083'
Set-TextValue $ws.Range("E18") '  -3.71%  '
Set-TextValue $ws.Range("D19") '0.06586'
Set-TextValue $ws.Range("E19") '  -0.84%  '
Set-TextValue $ws.Range("D20") '0.9984'
Set-TextValue $ws.Range("E20") '  -0.61%  '
Set-TextValue $ws.Range("D21") '17.34'
Set-TextValue $ws.Range("E21") '  -2.28%  '
Set-TextValue $ws.Range("D22") '6.028'
Set-TextValue $ws.Range("E22") '  -0.42%  '
Set-TextValue $ws.Range("D23") '28.345.48'
Set-TextValue $ws.Range("E23") '  -0.71%  '
Set-TextValue $ws.Range("D24") '11.19'
Set-TextValue $ws.Range("E24") '  -1.95%  '
Set-TextValue $ws.Range("D25") '2.227'
Set-TextValue $ws.Range("E25") '  -3.11%  '
Set-TextValue $ws.Range("D26") '160.88'
Set-TextValue $ws.Range("E26") '  +0.77%  '
Set-TextValue $ws.Range("D27") '20.62'
Set-TextValue $ws.Range("E27") '  -2.41%  '
Set-TextValue $ws.Range("D28") '2.011.40'
Set-TextValue $ws.Range("E28") '  -1.04%  '
Set-TextValue $ws.Range("D29") '2.407'
Set-TextValue $ws.Range("E29") '  +0.43%  '
Set-TextValue $ws.Range("D30") '128.34'
Set-TextValue $ws.Range("E30") '  +1.99%  '
Set-TextValue $ws.Range("D31") '0.1088'
Set-TextValue $ws.Range("E31") '  -0.23%  '
Set-TextValue $ws.Range("D32") '1.051'
Set-TextValue $ws.Range("E32") '  -4.24%  '
Set-TextValue $ws.Range("D33") '3.666'
Set-TextValue $ws.Range("E33") '  -0.37%  '
Set-TextValue $ws.Range("D34") '5.595'
Set-TextValue $ws.Range("E34") '  -2.27%  '
Set-TextValue $ws.Range("D35") '0.07168'
Set-TextValue $ws.Range("E35") '  -6.38%  '
Set-TextValue $ws.Range("D36") '9.109'
Set-TextValue $ws.Range("E36") '  +4.14%  '
Set-TextValue $ws.Range("D37") '0.02331'
Set-TextValue $ws.Range("E37") '  -1.85%  '
Set-TextValue $ws.Range("D38") '0.2160'
Set-TextValue $ws.Range("E38") '  -2.99%  '
Set-TextValue $ws.Range("E39") '  +1.40%  '
Set-TextValue $ws.Range("D40") '5.065'
Set-TextValue $ws.Range("E40") '  -3.89%  '
Set-TextValue $ws.Range("D41") '0.6239'
Set-TextValue $ws.Range("E41") '  -2.36%  '
Set-TextValue $ws.Range("D42") '0.9977'
Set-TextValue $ws.Range("E42") '  -0.54%  '
Set-TextValue $ws.Range("E43") '  -3.46%  '
Set-TextValue $ws.Range("D44") '1.325'
Set-TextValue $ws.Range("E44") '  -5.29%  '
Set-TextValue $ws.Range("D45") '13.18'
Set-TextValue $ws.Range("E45") '  -2.77%  '
Set-TextValue $ws.Range("D46") '0.5986'
Set-TextValue $ws.Range("E46") '  -2.14%  '
Set-TextValue $ws.Range("D47") '3.752'
Set-TextValue $ws.Range("E47") '  -1.25%  '
Set-TextValue $ws.Range("D48") '125.84'
Set-TextValue $ws.Range("E48") '  -1.47%  '
Set-TextValue $ws.Range("D49") '1.213'
Set-TextValue $ws.Range("E49") '  +0.79%  '
Set-TextValue $ws.Range("D50") '1.943'
Set-TextValue $ws.Range("E50") '  -2.64%  '
Set-TextValue $ws.Range("D51") '0.06862'
Set-TextValue $ws.Range("E51") '  -1.89%  '
